$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-02-20"

$ws.Range("I1").Value = "2022 (through 02-20)"

$ws.Range("I3").Value = 98
$ws.Range("I14").Value = 257
